# ---------------------------------------------------------------------------
# QA pass: root docs sync, figure refresh, registry updates
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsSources = $wb.Worksheets.Item("Sources")
$wsClaims  = $wb.Worksheets.Item("Claims")
$wsFigures = $wb.Worksheets.Item("Figures")

# ---------------------------------------------------------------------------
# Helper: write a value that must remain literal TEXT even when it looks like
# a number or an ISO date (e.g. "01", "2026-02-08"). Forcing the cell number
# format to Text ("@") before assignment stops Excel's auto-detection from
# converting the string into a numeric/date value.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# Sources sheet - rename two registry entries
# ---------------------------------------------------------------------------
$wsSources.Range("B117").Value = "Animal Health VC/PE Portfolio Mapping"
$wsSources.Range("B122").Value = "Global Ecosystem Landscape Map"

# ---------------------------------------------------------------------------
# Sources sheet - clear stray empty original_url / date_published cells
# for rows 448-463 (author instruction: never delete rows, only the
# extraneous blank cells)
# ---------------------------------------------------------------------------
for ($r = 448; $r -le 463; $r++) {
    $wsSources.Range("E$r").ClearContents()
    $wsSources.Range("F$r").ClearContents()
}

# ---------------------------------------------------------------------------
# Claims sheet - append four new QA-pass claims (rows 93-96)
# ---------------------------------------------------------------------------

# Row 93 - C092
$wsClaims.Range("A93").Value = "C092"
Set-TextValue $wsClaims.Range("B93") "01"
$wsClaims.Range("C93").Value = "Nutraceutical was first defined in 1989 as a food (or part of a food) providing medical or health benefits including disease prevention/treatment."
$wsClaims.Range("D93").Value = "S403"
$wsClaims.Range("E93").Value = "P3"
$wsClaims.Range("F93").Value = "N"
$wsClaims.Range("G93").Value = "Y"
Set-TextValue $wsClaims.Range("H93") "2026-02-08"
$wsClaims.Range("J93").Value = "Added definition sentence in executive summary during QA flow pass."

# Row 94 - C093
$wsClaims.Range("A94").Value = "C093"
$wsClaims.Range("B94").Value = "I.3"
$wsClaims.Range("C94").Value = "Global pet supplements were estimated at USD 2.26B in 2024 and projected to grow at 5.9% CAGR (2025-2030)."
$wsClaims.Range("D94").Value = "S416"
$wsClaims.Range("E94").Value = "Meta description"
$wsClaims.Range("F94").Value = "N"
$wsClaims.Range("G94").Value = "Y"
Set-TextValue $wsClaims.Range("H94") "2026-02-08"
$wsClaims.Range("J94").Value = "Added explicit market-size/CAGR quantification in Part I segment discussion."

# Row 95 - C094
$wsClaims.Range("A95").Value = "C094"
$wsClaims.Range("B95").Value = "I.1"
$wsClaims.Range("C95").Value = "Regulatory comparison figure expanded to four jurisdictions (US, EU, UK, China) with route-specific claim/entry differences."
$wsClaims.Range("D95").Value = "S085, S015, S124"
$wsClaims.Range("E95").Value = "Table / figure synthesis"
$wsClaims.Range("F95").Value = "N"
$wsClaims.Range("G95").Value = "Y"
Set-TextValue $wsClaims.Range("H95") "2026-02-08"
$wsClaims.Range("J95").Value = "Figure 1 scope widened and narrative updated in Part I."

# Row 96 - C095
$wsClaims.Range("A96").Value = "C095"
$wsClaims.Range("B96").Value = "I.3"
$wsClaims.Range("C96").Value = "Functional segment values in Part I include Gut Health USD 2,913M, Delivery Systems USD 2,749M, Immunity USD 1,841M, and Performance/FCR USD 1,426M."
$wsClaims.Range("D96").Value = "S089"
$wsClaims.Range("E96").Value = "Tab: Figure 18"
$wsClaims.Range("F96").Value = "N"
$wsClaims.Range("G96").Value = "Y"
Set-TextValue $wsClaims.Range("H96") "2026-02-08"
$wsClaims.Range("J96").Value = "Added during segment quantification enrichment pass."

# ---------------------------------------------------------------------------
# Figures sheet - refresh source mappings / notes on existing figures
# ---------------------------------------------------------------------------

# Row 2 - FIG-01
$wsFigures.Range("D2").Value = "S109, S110"
$wsFigures.Range("I2").Value = "Developed-market ownership chart regenerated; Mexico removed per QA."

# Row 21 - FIG-21
$wsFigures.Range("D21").Value = "S115, S116"
$wsFigures.Range("I21").Value = "Axis labeling clarified for pharma integration funnel."

# Row 23 - FIG-Figure_II_0_1_Innovation_Matrix
$wsFigures.Range("D23").Value = "S089, S115, S117, S118, S119, S120, S125"
$wsFigures.Range("I23").Value = "Regenerated with 15-company labeling and corrected high-R&D/high-margin quadrant cue."

# Row 38 - FIG-Table_US_vs_EU
$wsFigures.Range("D38").Value = "S085, S015, S124"
$wsFigures.Range("I38").Value = "Expanded to include UK and China in regulatory comparison."

# Row 46 - FIG-ES-1
$wsFigures.Range("D46").Value = "S121, S116"
$wsFigures.Range("I46").Value = "Source naming refined for clarity in executive summary."

# ---------------------------------------------------------------------------
# Figures sheet - two new simplified-index figures (rows 47-48)
# ---------------------------------------------------------------------------

# Row 47 - FIG-II-1-S
$wsFigures.Range("A47").Value = "FIG-II-1-S"
$wsFigures.Range("B47").Value = "Figure II.1 (Simplified): Functional needs across companion and production species."
$wsFigures.Range("C47").Value = "heatmap"
$wsFigures.Range("D47").Value = "S089, S116"
$wsFigures.Range("E47").Value = "Figure 18"
$wsFigures.Range("F47").Value = "Part I"
$wsFigures.Range("G47").Value = "active"
Set-TextValue $wsFigures.Range("H47") "2026-02-08"
$wsFigures.Range("I47").Value = "Added for simplified index readability per author QA pass."

# Row 48 - FIG-II-12-S
$wsFigures.Range("A48").Value = "FIG-II-12-S"
$wsFigures.Range("B48").Value = "Figure II.12 (Simplified): Comparative economic value by segment and primary sector."
$wsFigures.Range("C48").Value = "stacked_bar"
$wsFigures.Range("D48").Value = "S089"
$wsFigures.Range("E48").Value = "Figure 18"
$wsFigures.Range("F48").Value = "Part I"
$wsFigures.Range("G48").Value = "active"
Set-TextValue $wsFigures.Range("H48") "2026-02-08"
$wsFigures.Range("I48").Value = "Added for simplified index readability per author QA pass."
